$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price values remain plain text (avoid Excel auto-converting
# dotted numeric-looking strings like "1.002" into real numbers).
foreach ($addr in @("D2","D3","D5","D6","D8","D9","D10","D12","D13","D15","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D40","D41","D44","D45","D47","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.952.46'
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").Value = '1.648.05'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '309.87'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").Value = '  -2.24%  '
$ws.Range("D8").Value = '0.3801'
$ws.Range("E8").Value = '  -3.03%  '
$ws.Range("D9").Value = '51.91'
$ws.Range("E9").Value = '  -0.34%  '
$ws.Range("D10").Value = '1.341'
$ws.Range("E10").Value = '  -4.37%  '
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '0.08433'
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("D13").Value = '23.86'
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("E14").Value = '  -4.04%  '
$ws.Range("D15").Value = '8.016'
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").Value = '0.00001305'
$ws.Range("E16").Value = '  -4.01%  '
$ws.Range("D17").Value = '1.651.00'
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").Value = '94.04'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").Value = '0.06984'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '19.61'
$ws.Range("E20").Value = '  -4.82%  '
$ws.Range("D21").Value = '6.931'
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '13.72'
$ws.Range("D24").Value = '23.941.67'
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").Value = '2.454'
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").Value = '2.946'
$ws.Range("E26").Value = '  -3.10%  '
$ws.Range("E27").Value = '  -2.36%  '
$ws.Range("D28").Value = '153.45'
$ws.Range("E28").Value = '  -2.72%  '
$ws.Range("D29").Value = '5.396'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").Value = '137.95'
$ws.Range("E30").Value = '  -3.41%  '
$ws.Range("D31").Value = '7.823'
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").Value = '2.508'
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("D33").Value = '1.828.72'
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '1.015'
$ws.Range("E34").Value = '  -5.29%  '
$ws.Range("D35").Value = '0.08065'
$ws.Range("E35").Value = '  -2.79%  '
$ws.Range("D36").Value = '6.735'
$ws.Range("E36").Value = '  -2.61%  '
$ws.Range("D37").Value = '0.02928'
$ws.Range("E37").Value = '  -3.52%  '
$ws.Range("E38").Value = '  -3.86%  '
$ws.Range("E39").Value = '  -4.07%  '
$ws.Range("D40").Value = '0.09058'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("D41").Value = '0.7568'
$ws.Range("E41").Value = '  -2.35%  '
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").Value = '16.20'
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("D45").Value = '0.6937'
$ws.Range("E45").Value = '  -2.79%  '
$ws.Range("E46").Value = '  -4.01%  '
$ws.Range("D47").Value = '4.092'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").Value = '0.08293'
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("D50").Value = '133.92'
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("E51").Value = '  -3.91%  '
